# Fix create_current & create_legacy script names:
# "..._data_files_with_filename_column.sh" -> "..._data_file_with_filename_column.sh"
# for the specific scripts that had the typo (indi/outc/rpsr/ther under "current",
# and indi/outc/reac/rpsr/ther under "legacy").

$d = $word.ActiveDocument

$replacements = @(
    "create_current_all_indi_data_files_with_filename_column.sh",
    "create_current_all_outc_data_files_with_filename_column.sh",
    "create_current_all_rpsr_data_files_with_filename_column.sh",
    "create_current_all_ther_data_files_with_filename_column.sh",
    "create_legacy_all_indi_data_files_with_filename_column.sh",
    "create_legacy_all_outc_data_files_with_filename_column.sh",
    "create_legacy_all_reac_data_files_with_filename_column.sh",
    "create_legacy_all_rpsr_data_files_with_filename_column.sh",
    "create_legacy_all_ther_data_files_with_filename_column.sh"
)

foreach ($old in $replacements) {
    $new = $old -replace "_data_files_with_filename_column\.sh$", "_data_file_with_filename_column.sh"
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
